$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Se agrego la columna entorno" -- the user (email) values now carry an
# environment suffix (.uat / .int) to distinguish the target environment.
$ws.Range("F2").Value = "alex@bcp.com.pe.uat"
$ws.Range("F3").Value = "alex2@bcp.com.pe.uat"
$ws.Range("F4").Value = "alex3@bcp.com.pe.int"
$ws.Range("F5").Value = "juan1@bcp.com.pe.uat"
$ws.Range("F6").Value = "juan1@bcp.com.pe.int"

# Rebuild the hyperlinks (still pointing at the original mailto addresses)
# and restore the hyperlink cell style that gets reset on Add().
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:alex@bcp.com.pe")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:alex2@bcp.com.pe")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:alex3@bcp.com.pe")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:juan1@bcp.com.pe")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:juan1@bcp.com.pe")
$ws.Range("F2:F6").Style = "Hipervínculo"

# Move the active selection to F7
$ws.Range("F7").Select()
